# Update the "cryptos" price/volume snapshot to the latest scraped values.
# Note: several "Price" values look numeric (e.g. "339.08", "1.000",
# "0.00001044") but must stay exact text, matching the source data
# (trailing zeros, fixed decimal places, no scientific notation).
# A leading apostrophe forces Excel to store them as text instead of
# silently re-parsing/re-formatting them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.405.52'
$ws.Range("E2").Value = '  +3.39%  '
$ws.Range("D3").Value = '1.870.11'
$ws.Range("E3").Value = '  +2.04%  '
$ws.Range("E4").Value = '  -0.32%  '
$ws.Range("D5").Value = '''339.08'
$ws.Range("E5").Value = '  +2.25%  '
$ws.Range("E6").Value = '  -0.40%  '
$ws.Range("D7").Value = '''0.4700'
$ws.Range("E7").Value = '  +2.27%  '
$ws.Range("D8").Value = '''0.3969'
$ws.Range("E8").Value = '  +3.94%  '
$ws.Range("D9").Value = '''47.74'
$ws.Range("E9").Value = '  +2.44%  '
$ws.Range("D10").Value = '''0.08033'
$ws.Range("E10").Value = '  +1.50%  '
$ws.Range("D11").Value = '''1.001'
$ws.Range("E11").Value = '  +3.05%  '
$ws.Range("D12").Value = '''21.98'
$ws.Range("E12").Value = '  +4.31%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '''6.043'
$ws.Range("E13").Value = '  +2.56%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.870.07'
$ws.Range("E14").Value = '  +2.23%  '
$ws.Range("D15").Value = '''7.263'
$ws.Range("E15").Value = '  +3.34%  '
$ws.Range("D16").Value = '''91.14'
$ws.Range("E16").Value = '  +3.58%  '
$ws.Range("D17").Value = '''1.002'
$ws.Range("E17").Value = '  -0.39%  '
$ws.Range("D18").Value = '''0.00001044'
$ws.Range("E18").Value = '  +1.43%  '
$ws.Range("D19").Value = '''0.06627'
$ws.Range("E19").Value = '  -0.38%  '
$ws.Range("D20").Value = '''17.57'
$ws.Range("E20").Value = '  +3.34%  '
$ws.Range("E21").Value = '  -0.10%  '
$ws.Range("D22").Value = '28.404.09'
$ws.Range("E22").Value = '  +3.42%  '
$ws.Range("D23").Value = '''5.476'
$ws.Range("E23").Value = '  +2.40%  '
$ws.Range("D24").Value = '''11.07'
$ws.Range("E24").Value = '  +2.31%  '
$ws.Range("D25").Value = '''2.255'
$ws.Range("E25").Value = '  -2.32%  '
$ws.Range("D26").Value = '2.084.85'
$ws.Range("E26").Value = '  +1.70%  '
$ws.Range("D27").Value = '''160.72'
$ws.Range("E27").Value = '  +2.09%  '
$ws.Range("D28").Value = '''19.80'
$ws.Range("E28").Value = '  +2.22%  '
$ws.Range("D29").Value = '''2.125'
$ws.Range("E29").Value = '  +3.21%  '
$ws.Range("D30").Value = '''5.506'
$ws.Range("E30").Value = '  +3.93%  '
$ws.Range("D31").Value = '''120.32'
$ws.Range("E31").Value = '  +1.09%  '
$ws.Range("D32").Value = '''0.9737'
$ws.Range("E32").Value = '  +2.06%  '
$ws.Range("D33").Value = '''0.09512'
$ws.Range("E33").Value = '  +2.32%  '
$ws.Range("D34").Value = '''3.594'
$ws.Range("E34").Value = '  +0.10%  '
$ws.Range("D35").Value = '''1.376'
$ws.Range("E35").Value = '  +4.90%  '
$ws.Range("D36").Value = '''5.352'
$ws.Range("E36").Value = '  +2.02%  '
$ws.Range("D37").Value = '''0.06105'
$ws.Range("E37").Value = '  +2.90%  '
$ws.Range("D38").Value = '''0.02260'
$ws.Range("E38").Value = '  +3.08%  '
$ws.Range("D39").Value = '''8.367'
$ws.Range("E39").Value = '  +3.78%  '
$ws.Range("D40").Value = '''1.177'
$ws.Range("E40").Value = '  +1.22%  '
$ws.Range("D41").Value = '''0.5957'
$ws.Range("E41").Value = '  +2.92%  '
$ws.Range("D42").Value = '''1.000'
$ws.Range("E42").Value = '  -0.33%  '
$ws.Range("D43").Value = '''0.1881'
$ws.Range("E43").Value = '  +2.33%  '
$ws.Range("D44").Value = '''10.37'
$ws.Range("E44").Value = '  +3.43%  '
$ws.Range("D45").Value = '''1.285'
$ws.Range("E45").Value = '  +2.64%  '
$ws.Range("D46").Value = '''0.5591'
$ws.Range("D47").Value = '''12.09'
$ws.Range("E47").Value = '  +1.02%  '
$ws.Range("D48").Value = '''1.955'
$ws.Range("E48").Value = '  +4.94%  '
$ws.Range("D49").Value = '''0.06921'
$ws.Range("E49").Value = '  +4.17%  '
$ws.Range("D50").Value = '''2.061'
$ws.Range("E50").Value = '  +16.72%  '
$ws.Range("D51").Value = '''111.74'
$ws.Range("E51").Value = '  +1.21%  '
